# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
# Swap the match-data columns (B:AB) between specific row pairs, leaving
# column A (the running index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(9, 10),
    @(16, 17),
    @(92, 93),
    @(100, 101)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B${r1}:AB${r1}")
    $range2 = $ws.Range("B${r2}:AB${r2}")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value = $values2
    $range2.Value = $values1
}
